# "Math and MathTest ready"
# - New student name
# - Updated "don't change" note (now mentions indicating your name)
# - Rebalanced Final mark formula (Lab 4 & Lab 7 now double-weighted, /10 instead of /8)
# - Reworded / reordered "Important" notes block, new explanatory lines about lab scoring
# - Selection moved to C20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Student name -----------------------------------------------------
$ws.Range("A4").Value = "García Díaz, Vicente"

# --- Top note -----------------------------------------------------------
$ws.Range("A1").Value = "Please, don't change this file after you indicate your name"

# --- Final mark formula (single cell + table calculated column) --------
$ws.Range("J4").Formula = "=(B4+C4+D4+E4+F4+2*G4+H4+2*I4)/10"

# --- Reorganize the "Important:" notes block ----------------------------
# Move the existing explanatory sentences down so the two brand-new lines
# can be inserted above them (A15 "Maximum possible..." is dropped).
# A17 ("...<5 OR...") -> A18, then A16 ("...>=5 AND...") -> A17.
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4104) | Out-Null

$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4104) | Out-Null

$excel.CutCopyMode = 0

# New line 15: lab scoring explanation
$ws.Range("A15").Value = "Each lab will be graded with a score between 0 and 7 "

# New line 16: rich text, "Final mark" portion in bold
$ws.Range("A16").Value = "Total maximum possible Final mark without taking the final exam => 7 points"
$ws.Range("A16").Characters(24, 10).Font.Bold = $true
$ws.Range("A16").Characters(24, 10).Font.ColorIndex = -4105
$ws.Range("A16").Characters(34, 45).Font.Bold = $false
$ws.Range("A16").Characters(34, 45).Font.ColorIndex = -4105

# --- Selection ------------------------------------------------------------
$ws.Range("C20").Select()
